$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the latest cryptos snapshot: updated prices/volumes, and
# row 37/38 content swap (ARBITRUM <-> RenderToken).

$ws.Range("D2").Value = "46.313.29"
$ws.Range("E2").Value = "  +4.09%  "

$ws.Range("D3").Value = "2.457.35"
$ws.Range("E3").Value = "  +1.41%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.70"
$ws.Range("E5").Value = "  +2.59%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.44"
$ws.Range("E6").Value = "  +4.39%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.518"
$ws.Range("E7").Value = "  +0.97%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("E9").Value = "  +4.24%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.06"
$ws.Range("E10").Value = "  +2.26%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0809"
$ws.Range("E11").Value = "  +1.28%  "

$ws.Range("E12").Value = "  +0.66%  "

$ws.Range("E13").Value = "  -4.59%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.09"
$ws.Range("E14").Value = "  +2.30%  "

$ws.Range("D15").Value = "2.840.38"
$ws.Range("E15").Value = "  +1.43%  "

$ws.Range("D16").Value = "2.437.18"
$ws.Range("E16").Value = "  -1.02%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.845"
$ws.Range("E17").Value = "  +1.42%  "

$ws.Range("D18").Value = "46.143.51"
$ws.Range("E18").Value = "  +4.04%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.65"
$ws.Range("E19").Value = "  +1.77%  "

$ws.Range("E20").Value = "  +0.64%  "

$ws.Range("D21").Value = "0.0₃0937"
$ws.Range("E21").Value = "  +2.01%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.03"
$ws.Range("E22").Value = "  +3.36%  "

$ws.Range("E23").Value = "  +4.53%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "248.09"
$ws.Range("E24").Value = "  +2.77%  "

$ws.Range("E25").Value = "  +1.99%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.02"
$ws.Range("E26").Value = "  +3.47%  "

$ws.Range("E27").Value = "  -0.03%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.19"
$ws.Range("E28").Value = "  -4.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.71"
$ws.Range("E29").Value = "  +1.34%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.46"
$ws.Range("E30").Value = "  +3.66%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "49.38"
$ws.Range("E31").Value = "  +1.91%  "

$ws.Range("E32").Value = "  +5.93%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.89"
$ws.Range("E33").Value = "  +3.02%  "

$ws.Range("E34").Value = "  +3.38%  "

$ws.Range("E35").Value = "  -0.04%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0766"
$ws.Range("E36").Value = "  -0.19%  "

$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.91"
$ws.Range("E37").Value = "  +1.23%  "

$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.56"
$ws.Range("E38").Value = "  +1.98%  "

$ws.Range("E39").Value = "  +2.93%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "126.33"
$ws.Range("E40").Value = "  +4.46%  "

$ws.Range("E41").Value = "  +1.81%  "

$ws.Range("E42").Value = "  -0.05%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.85"
$ws.Range("E43").Value = "  +0.83%  "

$ws.Range("E44").Value = "  +1.26%  "

$ws.Range("D45").Value = "1.973.55"
$ws.Range("E45").Value = "  +1.53%  "

$ws.Range("E46").Value = "  +1.63%  "

$ws.Range("E47").Value = "  -4.30%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.87"
$ws.Range("E48").Value = "  +12.87%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.13"
$ws.Range("E49").Value = "  -3.58%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.07"
$ws.Range("E50").Value = "  +9.66%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.83"
$ws.Range("E51").Value = "  +6.13%  "
